# Update mods data [2026-02-07 15:13:58]
# Appends a new daily data row to the ModCounts sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModCounts")

$prevRow = 88
$newRow  = 89

# Leading apostrophe forces the date-like string to be stored as literal
# text (inlineStr) instead of being auto-converted to a date serial value.
$ws.Cells.Item($newRow, 1).Value = "'2026/02/07"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1181

# Copy the formatting (style) of the previous row onto the new row so the
# new cells share the same cell style as the rest of the data rows.
$src = $ws.Range("A$prevRow`:C$prevRow")
$dst = $ws.Range("A$newRow`:C$newRow")
$src.Copy() | Out-Null
$dst.PasteSpecial(-4122) | Out-Null
